# Update cryptocurrency price/volume snapshot data
# Commit: "Updated symbol list on Thu Jan 12 07:30:22 UTC 2023 with GitHub Actions"
#
# D (Price) and E (Volume 1h %) columns hold numeric-looking text that must
# stay plain text (as in the source inlineStr cells), so those values are
# written with a leading apostrophe (Excel text-prefix) to block automatic
# number/percentage conversion. B (Coin) / C (Link) columns are plain text
# already, so they are written as-is.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$textUpdates = @{
    "D2" = "285.46"
    "E2" = "3.03%"
    "D3" = "28.58"
    "E3" = "4.72%"
    "D4" = "4.996"
    "E4" = "2.62%"
    "D5" = "0.06489"
    "D6" = "7.227"
    "E6" = "4.17%"
    "D7" = "1.345"
    "E7" = "13.65%"
    "D8" = "0.9128"
    "E8" = "4.27%"
    "E9" = "-0.14%"
    "D10" = "0.06263"
    "E10" = "21.17%"
    "D11" = "0.07643"
    "E11" = "1.92%"
    "D12" = "0.02981"
    "E12" = "0.90%"
    "D13" = "0.08966"
    "E13" = "-0.13%"
    "D14" = "0.001592"
    "E14" = "1.65%"
    "D15" = "0.0006555"
    "E15" = "2.69%"
    "D16" = "0.006035"
    "E16" = "-1.39%"
    "D17" = "3.462"
    "E17" = "-0.47%"
    "E18" = "1.83%"
    "D19" = "2.242"
    "E19" = "-1.45%"
    "D20" = "0.3150"
    "E20" = "0.44%"
    "D21" = "0.1342"
    "E21" = "1.56%"
    "D22" = "3.982"
    "E22" = "2.06%"
    "D23" = "0.1554"
    "E23" = "3.58%"
    "D24" = "0.04455"
    "E24" = "0.95%"
    "D25" = "0.001189"
    "E25" = "1.08%"
    "D26" = "0.004330"
    "E26" = "11.99%"
    "D28" = "0.0001180"
    "E28" = "-9.24%"
    "E29" = "-15.76%"
    "D40" = "0.04144"
    "E40" = "-0.45%"
    "D41" = "0.006728"
    "E41" = "-0.98%"
    "D42" = "0.1232"
    "E42" = "5.01%"
    "D43" = "0.002139"
    "E43" = "10.31%"
    "D44" = "0.01176"
    "E44" = "2.60%"
    "D45" = "0.00005382"
    "E45" = "1.26%"
    "D46" = "0.01851"
    "E46" = "-0.04%"
    "D47" = "2.041"
    "E47" = "20.97%"
}

$plainUpdates = @{
    "B46" = "CoinbaseStockToken"
    "C46" = "https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin"
    "B47" = "BOLO"
    "C47" = "https://coinranking.com/coin/ogrGe0dEab+bolo-bolo"
}

foreach ($cellRef in $textUpdates.Keys) {
    $ws.Range($cellRef).Value = "'" + $textUpdates[$cellRef]
}

foreach ($cellRef in $plainUpdates.Keys) {
    $ws.Range($cellRef).Value = $plainUpdates[$cellRef]
}

